$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update cryptocurrency price/volume data per latest scrape
# Leading apostrophe forces Excel to store values as literal text
# (matches original inlineStr cell type, avoids numeric auto-conversion)

# Row 2
$ws.Range("D2").Value = "'38.018.70"
$ws.Range("E2").Value = "'  +2.59%  "

# Row 3
$ws.Range("D3").Value = "'2.056.11"
$ws.Range("E3").Value = "'  +1.93%  "

# Row 4
$ws.Range("E4").Value = "'  +0.15%  "

# Row 5
$ws.Range("D5").Value = "'230.11"
$ws.Range("E5").Value = "'  +1.77%  "

# Row 6
$ws.Range("E6").Value = "'  +2.01%  "

# Row 7
$ws.Range("D7").Value = "'58.21"
$ws.Range("E7").Value = "'  +6.33%  "

# Row 8
$ws.Range("E8").Value = "'  +0.02%  "

# Row 9
$ws.Range("D9").Value = "'0.386"
$ws.Range("E9").Value = "'  +2.86%  "

# Row 10
$ws.Range("E10").Value = "'  +3.04%  "

# Row 11
$ws.Range("E11").Value = "'  +1.19%  "

# Row 12
$ws.Range("D12").Value = "'2.360.01"
$ws.Range("E12").Value = "'  +1.98%  "

# Row 13
$ws.Range("D13").Value = "'14.63"
$ws.Range("E13").Value = "'  +3.68%  "

# Row 14
$ws.Range("D14").Value = "'20.68"
$ws.Range("E14").Value = "'  +2.43%  "

# Row 15
$ws.Range("E15").Value = "'  +2.15%  "

# Row 16
$ws.Range("D16").Value = "'5.27"
$ws.Range("E16").Value = "'  +2.80%  "

# Row 17
$ws.Range("D17").Value = "'2.063.20"
$ws.Range("E17").Value = "'  +2.11%  "

# Row 18
$ws.Range("D18").Value = "'37.935.81"
$ws.Range("E18").Value = "'  +2.54%  "

# Row 19
$ws.Range("D19").Value = "'6.14"
$ws.Range("E19").Value = "'  -1.08%  "

# Row 20
$ws.Range("D20").Value = "'69.76"
$ws.Range("E20").Value = "'  +1.26%  "

# Row 21
$ws.Range("D21").Value = "'0.0₃0831"
$ws.Range("E21").Value = "'  +1.87%  "

# Row 22
$ws.Range("D22").Value = "'224.90"
$ws.Range("E22").Value = "'  +0.84%  "

# Row 23
$ws.Range("E23").Value = "'  +0.03%  "

# Row 24
$ws.Range("D24").Value = "'2.45"
$ws.Range("E24").Value = "'  +0.51%  "

# Row 25
$ws.Range("E25").Value = "'  +3.10%  "

# Row 26
$ws.Range("D26").Value = "'9.30"
$ws.Range("E26").Value = "'  +1.58%  "

# Row 27
$ws.Range("D27").Value = "'166.20"
$ws.Range("E27").Value = "'  +0.03%  "

# Row 28
$ws.Range("E28").Value = "'  +8.02%  "

# Row 29
$ws.Range("D29").Value = "'19.04"
$ws.Range("E29").Value = "'  +1.82%  "

# Row 30
$ws.Range("E30").Value = "'  +0.44%  "

# Row 31
$ws.Range("E31").Value = "'  +1.87%  "

# Row 32
$ws.Range("D32").Value = "'4.54"
$ws.Range("E32").Value = "'  +0.87%  "

# Row 33
$ws.Range("E33").Value = "'  +3.85%  "

# Row 34
$ws.Range("D34").Value = "'0.0612"
$ws.Range("E34").Value = "'  +0.11%  "

# Row 35
$ws.Range("E35").Value = "'  +9.23%  "

# Row 36
$ws.Range("D36").Value = "'2.34"
$ws.Range("E36").Value = "'  +0.08%  "

# Row 37
$ws.Range("D37").Value = "'5.98"
$ws.Range("E37").Value = "'  +13.11%  "

# Row 38
$ws.Range("E38").Value = "'  +5.40%  "

# Row 39
$ws.Range("E39").Value = "'  +0.12%  "

# Row 40
$ws.Range("D40").Value = "'98.45"
$ws.Range("E40").Value = "'  +3.83%  "

# Row 41
$ws.Range("E41").Value = "'  +1.56%  "

# Row 42
$ws.Range("D42").Value = "'1.486.05"

# Row 43
$ws.Range("E43").Value = "'  +3.31%  "

# Row 44
$ws.Range("D44").Value = "'0.0937"
$ws.Range("E44").Value = "'  +2.75%  "

# Row 45
$ws.Range("D45").Value = "'16.66"
$ws.Range("E45").Value = "'  +2.67%  "

# Row 46
$ws.Range("B46").Value = "'FTXToken"
$ws.Range("C46").Value = "'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
$ws.Range("D46").Value = "'4.17"
$ws.Range("E46").Value = "'  +19.34%  "

# Row 47
$ws.Range("B47").Value = "'TrustWalletToken"
$ws.Range("C47").Value = "'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D47").Value = "'1.13"
$ws.Range("E47").Value = "'  +0.80%  "

# Row 48
$ws.Range("E48").Value = "'  +0.87%  "

# Row 49
$ws.Range("E49").Value = "'  +1.72%  "

# Row 50
$ws.Range("E50").Value = "'  -1.84%  "

# Row 51
$ws.Range("D51").Value = "'2.247.92"
$ws.Range("E51").Value = "'  +2.22%  "

